$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header "OMIN-Gene" -> "OMIM"
$ws.Range("J1").Value = "OMIM"

# Add new header column "KEGG" in R1, matching style of neighboring headers
$ws.Range("P1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "KEGG"

# Match the author's saved view/selection state
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("R2").Select()
